$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $cell = $ws.Range("A1")
    if ($cell.Value2 -like "Month/Year:*") {
        $cell.Value = "Month/Year: FEBRUARY 2023"
        $ws.Rows.Item(1).AutoFit()
    }
}
